# Updates the cryptos list (prices / 1h volume %) to the latest snapshot.
# Column D values that look numeric are entered with a leading apostrophe
# so Excel stores them as text (matching the original inline-string cells)
# instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.652.64'
$ws.Range("E2").Value = '  +1.62%  '

$ws.Range("D3").Value = '2.296.71'
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''323.54'
$ws.Range("E5").Value = '  +1.93%  '

$ws.Range("D6").Value = '''104.52'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").Value = '''0.611'
$ws.Range("E9").Value = '  +0.50%  '

$ws.Range("D10").Value = '''40.26'
$ws.Range("E10").Value = '  +2.14%  '

$ws.Range("D11").Value = '''0.0909'
$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").Value = '''8.55'
$ws.Range("E12").Value = '  +2.01%  '

$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = '''0.972'
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").Value = '''15.27'
$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").Value = '2.648.16'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '2.295.21'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").Value = '42.621.56'
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = '''7.46'
$ws.Range("E19").Value = '  -2.65%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").Value = '''13.51'
$ws.Range("E21").Value = '  +34.71%  '

$ws.Range("D22").Value = '''73.34'
$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").Value = '''3.59'
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").Value = '''269.39'
$ws.Range("E24").Value = '  -6.70%  '

$ws.Range("D25").Value = '''2.23'
$ws.Range("E25").Value = '  -2.01%  '

$ws.Range("E26").Value = '  -0.25%  '

$ws.Range("D27").Value = '''10.92'
$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").Value = '''2.33'
$ws.Range("E28").Value = '  +2.97%  '

$ws.Range("D29").Value = '''38.73'
$ws.Range("E29").Value = '  +9.83%  '

$ws.Range("D30").Value = '''22.57'
$ws.Range("E30").Value = '  -3.63%  '

$ws.Range("D31").Value = '''165.62'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("E32").Value = '  +5.25%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("E34").Value = '  +0.70%  '

$ws.Range("E35").Value = '  -1.70%  '

$ws.Range("E36").Value = '  -13.74%  '

$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("E38").Value = '  +1.63%  '

$ws.Range("D39").Value = '''3.75'
$ws.Range("E39").Value = '  +3.79%  '

$ws.Range("D40").Value = '''2.75'
$ws.Range("E40").Value = '  -3.49%  '

$ws.Range("E41").Value = '  +4.44%  '

$ws.Range("D42").Value = '''70.04'
$ws.Range("E42").Value = '  -1.00%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.226'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D45").Value = '''93.26'
$ws.Range("E45").Value = '  -9.86%  '

$ws.Range("D46").Value = '''12.34'
$ws.Range("E46").Value = '  +2.24%  '

$ws.Range("D47").Value = '''81.68'
$ws.Range("E47").Value = '  +4.26%  '

$ws.Range("D48").Value = '''113.91'
$ws.Range("E48").Value = '  -2.11%  '

$ws.Range("D49").Value = '''8.94'
$ws.Range("E49").Value = '  -1.56%  '

$ws.Range("D50").Value = '''5.28'
$ws.Range("E50").Value = '  -1.05%  '

$ws.Range("D51").Value = '1.594.15'
$ws.Range("E51").Value = '  +2.29%  '
